# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.129.48'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +5.74%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.920.82'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +2.67%  '
$ws.Range("E4").Value = '  -0.54%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '330.44'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +4.79%  '
$ws.Range("E6").Value = '  -0.56%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5204'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +2.19%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4081'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +4.41%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08524'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +2.01%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '43.30'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +3.84%  '
$ws.Range("E11").Value = '  +2.25%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.38'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +9.85%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.415'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +3.23%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.920.36'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +2.34%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.412'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.84%  '
$ws.Range("E16").Value = '  -0.47%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '95.73'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +5.04%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001114'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.97%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06726'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.04%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.32'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +3.33%  '
$ws.Range("E21").Value = '  -0.52%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.031'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +2.13%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '30.140.15'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +5.67%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.37'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +2.11%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.222'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.11%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.142.59'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +2.54%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.14'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +2.45%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '159.99'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.67%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.457'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +2.16%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '129.05'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +1.58%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.078'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +3.84%  '
$ws.Range("E32").Value = '  +1.38%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.097'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +6.34%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.632'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.71%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02496'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +1.81%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06613'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.91%  '
$ws.Range("E37").Value = '  +2.41%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.218'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +4.02%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.236'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +4.98%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.941'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.04%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6534'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +2.49%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.69'
$ws.Range("D42").ClearFormats()
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.249'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +1.37%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6172'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +2.91%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.21'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +1.79%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.773'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +2.41%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.090'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +4.54%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.251'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +2.92%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '124.79'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +2.41%  '
$ws.Range("E50").Value = '  +2.71%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '79.72'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +4.40%  '
